# Append new trading-log rows (154-157) to Sheet1, mirroring the
# "real_trading_log.xlsx" export format: one row per logged event.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 154: TRADING_ATTEMPT for BTC
$ws.Cells.Item(154, 1).Value = "2026-01-05T01:55:52.741161"
$ws.Cells.Item(154, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(154, 3).Value = "BTC"
$ws.Cells.Item(154, 4).Value = "UNKNOWN"
$ws.Cells.Item(154, 5).Value = 93117.57494500774
$ws.Cells.Item(154, 6).Value = ""
$ws.Cells.Item(154, 7).Value = ""
$ws.Cells.Item(154, 8).Value = ""
$ws.Cells.Item(154, 9).Value = ""
$ws.Cells.Item(154, 10).Value = ""
$ws.Cells.Item(154, 11).Value = "ATTEMPT"
$ws.Cells.Item(154, 12).Value = "Attempting trade 1/2"

# Row 155: POSITION_OPENED for BTC
$ws.Cells.Item(155, 1).Value = "2026-01-05T01:55:54.403711"
$ws.Cells.Item(155, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(155, 3).Value = "BTC"
$ws.Cells.Item(155, 4).Value = "UNKNOWN"
$ws.Cells.Item(155, 5).Value = 93117.57494500774
$ws.Cells.Item(155, 6).Value = 7200
$ws.Cells.Item(155, 7).Value = 40
$ws.Cells.Item(155, 8).Value = 1.681936884743757
$ws.Cells.Item(155, 9).Value = ""
$ws.Cells.Item(155, 10).Value = ""
$ws.Cells.Item(155, 11).Value = "SUCCESS"
$ws.Cells.Item(155, 12).Value = ""

# Row 156: TRADING_ATTEMPT for ETH
$ws.Cells.Item(156, 1).Value = "2026-01-05T01:55:54.448482"
$ws.Cells.Item(156, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(156, 3).Value = "ETH"
$ws.Cells.Item(156, 4).Value = "UNKNOWN"
$ws.Cells.Item(156, 5).Value = 3204.677318622782
$ws.Cells.Item(156, 6).Value = ""
$ws.Cells.Item(156, 7).Value = ""
$ws.Cells.Item(156, 8).Value = ""
$ws.Cells.Item(156, 9).Value = ""
$ws.Cells.Item(156, 10).Value = ""
$ws.Cells.Item(156, 11).Value = "ATTEMPT"
$ws.Cells.Item(156, 12).Value = "Attempting trade 2/2"

# Row 157: POSITION_FAILED for ETH
$ws.Cells.Item(157, 1).Value = "2026-01-05T01:55:54.682973"
$ws.Cells.Item(157, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(157, 3).Value = "ETH"
$ws.Cells.Item(157, 4).Value = "UNKNOWN"
$ws.Cells.Item(157, 5).Value = ""
$ws.Cells.Item(157, 6).Value = ""
$ws.Cells.Item(157, 7).Value = ""
$ws.Cells.Item(157, 8).Value = ""
$ws.Cells.Item(157, 9).Value = ""
$ws.Cells.Item(157, 10).Value = ""
$ws.Cells.Item(157, 11).Value = "FAILED"
$ws.Cells.Item(157, 12).Value = "Trade execution failed for trade 2"
